$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(
    6.9774387942869378,
    6.004243779615428,
    7.6322864032111246,
    6.7936430562505983,
    10.211648421334095,
    10.036886602677358,
    5.192138874353029,
    9.0180021563608577,
    6.4006940188468775,
    7.2590030465059279,
    7.9571657832797493,
    10.161989665999648,
    7.2305734475742787,
    10.387100258666949
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
